# Updates cryptos list values (prices / 1h volume %) per the commit diff.
# Some D-column price strings look numeric (e.g. "1.00", "74.45") but must
# stay plain text (matching the source workbook, where every cell is an
# inline string) -- NumberFormat "@" is set first so Excel does not coerce
# them into numbers when .Value is assigned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '44.190.90'
$ws.Range("E2").Value = '  +1.56%  '

# Row 3
$ws.Range("D3").Value = '2.359.75'
$ws.Range("E3").Value = '  -0.92%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.683'
$ws.Range("E5").Value = '  +3.95%  '

# Row 6
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.15'
$ws.Range("E6").Value = '  +3.66%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.45'
$ws.Range("E7").Value = '  +3.27%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.583'
$ws.Range("E9").Value = '  +25.06%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.102'
$ws.Range("E10").Value = '  +5.01%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '31.92'
$ws.Range("E11").Value = '  +18.66%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.49'
$ws.Range("E12").Value = '  +19.65%  '

# Row 13
$ws.Range("E13").Value = '  +2.10%  '

# Row 14
$ws.Range("D14").Value = '2.709.71'
$ws.Range("E14").Value = '  -0.81%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.92'
$ws.Range("E15").Value = '  +5.79%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.915'
$ws.Range("E16").Value = '  +6.61%  '

# Row 17
$ws.Range("D17").Value = '2.352.68'
$ws.Range("E17").Value = '  -1.43%  '

# Row 18
$ws.Range("D18").Value = '44.395.58'
$ws.Range("E18").Value = '  +2.19%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000104'
$ws.Range("E19").Value = '  +4.41%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.76'
$ws.Range("E20").Value = '  +5.73%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '78.40'
$ws.Range("E21").Value = '  +5.12%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '256.18'
$ws.Range("E22").Value = '  +1.72%  '

# Row 23
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.03%  '

# Row 24
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.59'
$ws.Range("E24").Value = '  +4.48%  '

# Row 25
$ws.Range("E25").Value = '  -5.10%  '

# Row 26
$ws.Range("E26").Value = '  +7.30%  '

# Row 27
$ws.Range("E27").Value = '  +3.57%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.59'
$ws.Range("E28").Value = '  -2.55%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.62'
$ws.Range("E29").Value = '  +4.76%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.85'
$ws.Range("E30").Value = '  +0.12%  '

# Row 31
$ws.Range("E31").Value = '  +3.77%  '

# Row 32
$ws.Range("E32").Value = '  +4.76%  '

# Row 33
$ws.Range("E33").Value = '  +8.13%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0759'
$ws.Range("E34").Value = '  +9.26%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.37'
$ws.Range("E35").Value = '  +4.98%  '

# Row 36
$ws.Range("E36").Value = '  +6.60%  '

# Row 37
$ws.Range("E37").Value = '  +0.34%  '

# Row 38
$ws.Range("E38").Value = '  -0.41%  '

# Row 39
$ws.Range("E39").Value = '  +7.30%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.46'
$ws.Range("E40").Value = '  +4.08%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.02'
$ws.Range("E41").Value = '  +0.78%  '

# Row 42
$ws.Range("E42").Value = '  -0.01%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.194'
$ws.Range("E43").Value = '  +14.96%  '

# Row 44
$ws.Range("E44").Value = '  +2.86%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.51'
$ws.Range("E45").Value = '  +11.11%  '

# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.93'
$ws.Range("E46").Value = '  +1.92%  '

# Row 47
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0996'
$ws.Range("E47").Value = '  +4.63%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.17'
$ws.Range("E48").Value = '  -1.10%  '

# Row 49
$ws.Range("E49").Value = '  +0.50%  '

# Row 50
$ws.Range("D50").Value = '1.457.26'
$ws.Range("E50").Value = '  +0.16%  '

# Row 51
$ws.Range("E51").Value = '  +4.78%  '
